$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "70.965.09"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +6.13%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.665.31"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +18.41%  "

$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "623.60"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +8.16%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "180.83"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.44%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.662.20"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +18.41%  "

$ws.Range("E8").Value = "  +0.04%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.538"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.77%  "

$ws.Range("E10").Value = "  +7.71%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.65"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.18%  "

$ws.Range("E12").Value = "  +7.08%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "40.58"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +12.42%  "

$ws.Range("E14").Value = "  +5.61%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.277.24"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +18.49%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.670.25"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +18.76%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "70.972.83"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +6.27%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.124"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.55%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.54"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +7.66%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "521.83"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +8.40%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.93"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.01%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.26"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +19.50%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.743"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +7.53%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "88.44"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.86%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.48"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +10.39%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "13.44"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +6.06%  "

$ws.Range("E27").Value = "  +8.38%  "

$ws.Range("E28").Value = "  -0.06%  "

$ws.Range("E29").Value = "  +11.32%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.10"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.56%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.92"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +12.31%  "

$ws.Range("E32").Value = "  +12.98%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0000109"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +16.33%  "

$ws.Range("E34").Value = "  +3.16%  "

$ws.Range("E35").Value = "  +0.10%  "

$ws.Range("E36").Value = "  +9.47%  "

$ws.Range("E37").Value = "  +8.25%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.347"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +11.15%  "

$ws.Range("E39").Value = "  +11.42%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "51.88"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.88%  "

$ws.Range("E41").Value = "  +5.07%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "45.57"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.96%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.82"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.90%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.121.49"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +11.74%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "423.61"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +14.01%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.77"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.46%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "28.58"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +14.64%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0369"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +7.42%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "139.29"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.75%  "

$ws.Range("E51").Value = "  +9.59%  "
